$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.108
$ws.Range("A3").Value = -21.666
$ws.Range("C3").Value = -13.204
$ws.Range("D6").Value = -7.804
$ws.Range("E8").Value = 16.777
$ws.Range("C12").Value = -11.169
$ws.Range("A14").Value = -21.603
$ws.Range("D19").Value = -7.859999999999999
$ws.Range("A21").Value = -20.257
$ws.Range("A23").Value = -20.775
$ws.Range("E23").Value = 16.46
$ws.Range("C24").Value = -12.459
$ws.Range("D24").Value = -7.435
$ws.Range("A25").Value = -20.536
$ws.Range("B25").Value = 6.851000000000001
$ws.Range("C25").Value = -12.406
$ws.Range("A26").Value = -21.368
$ws.Range("E26").Value = 16.576
$ws.Range("B27").Value = 5.626
$ws.Range("A29").Value = -21.279
$ws.Range("D30").Value = -7.191
$ws.Range("B31").Value = 5.355
$ws.Range("D31").Value = -8.086
$ws.Range("D33").Value = -7.851999999999999
$ws.Range("E37").Value = 16.468
$ws.Range("B39").Value = 7.657999999999999
$ws.Range("D42").Value = -8.431000000000001
$ws.Range("B48").Value = 5.176
$ws.Range("E48").Value = 17.349
$ws.Range("C50").Value = -13.133
$ws.Range("B51").Value = 5.42
$ws.Range("B52").Value = 5.358000000000001
$ws.Range("A53").Value = -21.642
$ws.Range("C53").Value = -11.646
$ws.Range("B55").Value = 4.508
$ws.Range("D55").Value = -8.210000000000001
$ws.Range("B56").Value = 5.144
$ws.Range("A57").Value = -21.322
$ws.Range("B57").Value = 7.019
$ws.Range("C57").Value = -12.872
$ws.Range("D58").Value = -7.970999999999999
$ws.Range("A59").Value = -22.097
$ws.Range("C61").Value = -13.025
$ws.Range("E62").Value = 16.656
$ws.Range("C63").Value = -11.466
$ws.Range("D65").Value = -7.869
$ws.Range("E66").Value = 17.18
$ws.Range("A69").Value = -21.56
$ws.Range("C70").Value = -12.201
$ws.Range("D70").Value = -7.449000000000001
$ws.Range("B73").Value = 6.981
$ws.Range("D75").Value = -7.645
$ws.Range("A79").Value = -21.215
$ws.Range("A83").Value = -22.185
$ws.Range("D83").Value = -8.397
$ws.Range("C86").Value = -12.416
$ws.Range("D86").Value = -7.851000000000001
$ws.Range("B89").Value = 4.999999999999999
$ws.Range("E89").Value = 17.081
$ws.Range("B90").Value = 5.917
$ws.Range("A91").Value = -21.527
$ws.Range("B92").Value = 6.058999999999999
$ws.Range("A93").Value = -21.439
$ws.Range("E94").Value = 17.498
$ws.Range("D96").Value = -7.519
$ws.Range("D97").Value = -8.17
$ws.Range("C98").Value = -12.45
$ws.Range("C100").Value = -12.183
$ws.Range("C102").Value = -13.143
